$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.208.64'
$ws.Range('E2').Value = '  -1.92%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.580.04'
$ws.Range('E3').Value = '  -1.33%  '
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '209.16'
$ws.Range('E5').Value = '  -1.15%  '
$ws.Range('E6').Value = '  -2.42%  '
$ws.Range('E7').Value = '  -0.31%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.0610'
$ws.Range('E8').Value = '  -1.56%  '
$ws.Range('E9').Value = '  -0.60%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.48'
$ws.Range('E10').Value = '  -1.02%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0845'
$ws.Range('E11').Value = '  +0.01%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.802.98'
$ws.Range('E12').Value = '  -1.27%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.586.87'
$ws.Range('E13').Value = '  -0.93%  '
$ws.Range('E14').Value = '  -0.40%  '
$ws.Range('E15').Value = '  -1.34%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.41'
$ws.Range('E16').Value = '  -0.77%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.195.59'
$ws.Range('E17').Value = '  -1.86%  '
$ws.Range('E18').Value = '  -1.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.24'
$ws.Range('E20').Value = '  -0.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '206.30'
$ws.Range('E21').Value = '  -1.81%  '
$ws.Range('E22').Value = '  -0.83%  '
$ws.Range('E23').Value = '  -3.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.83'
$ws.Range('E24').Value = '  -1.26%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.87'
$ws.Range('E25').Value = '  +0.42%  '
$ws.Range('E26').Value = '  -0.38%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.00'
$ws.Range('E27').Value = '  -1.19%  '
$ws.Range('E28').Value = '  -1.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.20'
$ws.Range('E29').Value = '  -1.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0502'
$ws.Range('E30').Value = '  -1.69%  '
$ws.Range('E31').Value = '  -1.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.21'
$ws.Range('E32').Value = '  -1.67%  '
$ws.Range('E33').Value = '  -1.40%  '
$ws.Range('B34').Value = 'Maker'
$ws.Range('C34').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.281.62'
$ws.Range('E34').Value = '  -1.14%  '
$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.25'
$ws.Range('E35').Value = '  +7.60%  '
$ws.Range('E36').Value = '  -0.16%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.603'
$ws.Range('E37').Value = '  +0.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.48'
$ws.Range('E38').Value = '  -1.21%  '
$ws.Range('E39').Value = '  -2.18%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.812'
$ws.Range('E40').Value = '  -1.88%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.52'
$ws.Range('E41').Value = '  +2.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.768'
$ws.Range('E42').Value = '  -1.50%  '
$ws.Range('E43').Value = '  -3.00%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '62.00'
$ws.Range('E44').Value = '  -1.69%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.715.95'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '88.44'
$ws.Range('E46').Value = '  -2.45%  '
$ws.Range('E47').Value = '  -0.54%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.101'
$ws.Range('E48').Value = '  -0.70%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0507'
$ws.Range('E49').Value = '  -2.00%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0₇0958'
$ws.Range('E50').Value = '  -9.85%  '
$ws.Range('E51').Value = '  -0.01%  '
